# Updated cryptos list: apply new Price (D) / Volume(1h) (E) text values.
# Values that look like plain numbers are forced to stay text (matching the
# source data's inlineStr formatting) by briefly applying a text NumberFormat
# and then resetting the cell style so no residual number formatting sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: cell reference, new value
$updates = @(
    @{ Cell = 'D2'; Value = '62.913.76' }
    @{ Cell = 'E2'; Value = '  +2.62%  ' }
    @{ Cell = 'D3'; Value = '2.950.30' }
    @{ Cell = 'E3'; Value = '  +0.70%  ' }
    @{ Cell = 'D4'; Value = '0.999' }
    @{ Cell = 'E4'; Value = '  +0.02%  ' }
    @{ Cell = 'D5'; Value = '593.21' }
    @{ Cell = 'E5'; Value = '  -0.33%  ' }
    @{ Cell = 'D6'; Value = '147.58' }
    @{ Cell = 'E6'; Value = '  +2.90%  ' }
    @{ Cell = 'E7'; Value = '  +0.01%  ' }
    @{ Cell = 'D8'; Value = '2.946.88' }
    @{ Cell = 'E8'; Value = '  +0.64%  ' }
    @{ Cell = 'D9'; Value = '0.505' }
    @{ Cell = 'E9'; Value = '  +0.82%  ' }
    @{ Cell = 'D10'; Value = '7.06' }
    @{ Cell = 'E10'; Value = '  +1.57%  ' }
    @{ Cell = 'E11'; Value = '  +5.60%  ' }
    @{ Cell = 'E12'; Value = '  +0.56%  ' }
    @{ Cell = 'E13'; Value = '  +4.41%  ' }
    @{ Cell = 'D14'; Value = '32.65' }
    @{ Cell = 'E14'; Value = '  -1.78%  ' }
    @{ Cell = 'E15'; Value = '  -0.79%  ' }
    @{ Cell = 'D16'; Value = '3.440.82' }
    @{ Cell = 'D17'; Value = '62.878.97' }
    @{ Cell = 'E17'; Value = '  +2.63%  ' }
    @{ Cell = 'E18'; Value = '  +0.31%  ' }
    @{ Cell = 'D19'; Value = '2.952.01' }
    @{ Cell = 'E19'; Value = '  +0.75%  ' }
    @{ Cell = 'D20'; Value = '439.66' }
    @{ Cell = 'E20'; Value = '  +1.49%  ' }
    @{ Cell = 'D21'; Value = '13.44' }
    @{ Cell = 'E21'; Value = '  -0.51%  ' }
    @{ Cell = 'D22'; Value = '0.665' }
    @{ Cell = 'E22'; Value = '  -1.04%  ' }
    @{ Cell = 'E23'; Value = '  -0.62%  ' }
    @{ Cell = 'D24'; Value = '11.16' }
    @{ Cell = 'E24'; Value = '  +3.22%  ' }
    @{ Cell = 'D25'; Value = '80.79' }
    @{ Cell = 'E25'; Value = '  -0.91%  ' }
    @{ Cell = 'E26'; Value = '  -1.25%  ' }
    @{ Cell = 'D27'; Value = '11.83' }
    @{ Cell = 'E27'; Value = '  +0.91%  ' }
    @{ Cell = 'E28'; Value = '  +0.01%  ' }
    @{ Cell = 'E29'; Value = '  +1.68%  ' }
    @{ Cell = 'E30'; Value = '  +6.01%  ' }
    @{ Cell = 'D31'; Value = '2.61' }
    @{ Cell = 'E31'; Value = '  +0.72%  ' }
    @{ Cell = 'D32'; Value = '0.0000102' }
    @{ Cell = 'E32'; Value = '  +17.23%  ' }
    @{ Cell = 'D33'; Value = '26.36' }
    @{ Cell = 'E33'; Value = '  -0.90%  ' }
    @{ Cell = 'E34'; Value = '  -0.78%  ' }
    @{ Cell = 'D35'; Value = '0.999' }
    @{ Cell = 'E35'; Value = '  -0.01%  ' }
    @{ Cell = 'D36'; Value = '0.991' }
    @{ Cell = 'E36'; Value = '  -2.11%  ' }
    @{ Cell = 'E37'; Value = '  -0.33%  ' }
    @{ Cell = 'E38'; Value = '  +3.27%  ' }
    @{ Cell = 'D39'; Value = '49.70' }
    @{ Cell = 'E39'; Value = '  -0.12%  ' }
    @{ Cell = 'E40'; Value = '  +1.85%  ' }
    @{ Cell = 'D41'; Value = '8.45' }
    @{ Cell = 'E41'; Value = '  -0.57%  ' }
    @{ Cell = 'E42'; Value = '  -3.52%  ' }
    @{ Cell = 'D43'; Value = '0.279' }
    @{ Cell = 'E43'; Value = '  +0.19%  ' }
    @{ Cell = 'D44'; Value = '39.30' }
    @{ Cell = 'E44'; Value = '  -6.31%  ' }
    @{ Cell = 'D45'; Value = '2.707.05' }
    @{ Cell = 'E45'; Value = '  +0.52%  ' }
    @{ Cell = 'D46'; Value = '135.09' }
    @{ Cell = 'E46'; Value = '  +1.02%  ' }
    @{ Cell = 'D47'; Value = '0.0338' }
    @{ Cell = 'E47'; Value = '  -1.83%  ' }
    @{ Cell = 'D48'; Value = '359.64' }
    @{ Cell = 'E48'; Value = '  -0.50%  ' }
    @{ Cell = 'E50'; Value = '  -0.56%  ' }
    @{ Cell = 'D51'; Value = '22.73' }
    @{ Cell = 'E51'; Value = '  -3.29%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Value -match '^[+-]?\d+(\.\d+)?([eE][+-]?\d+)?$') {
        # Would otherwise be auto-converted to a number (and possibly
        # re-formatted, e.g. scientific notation) - keep it literal text.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}

Write-Output "Applied $($updates.Count) cell updates"
